$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the daily conversion note text in A1 ---
$hoja1 = $wb.Worksheets.Item("Hoja1")
$note = $hoja1.Range("A1").Value2
$note = $note.Replace("1000 Bs = 9.12 = 37828.47 pesos", "1000 Bs = 9.18 = 38148.42 pesos")
$note = $note.Replace("37828.47 pesos = 9.08 = 960.01 Bs", "38148.42 pesos = 9.17 = 975.1 Bs")
$hoja1.Range("A1").Value = $note

# --- Sheet "tasas": update the rate cells ---
$tasas = $wb.Worksheets.Item("tasas")
$tasas.Range("N10").Value = 108.88
$tasas.Range("O10").Value = 4153.6
$tasas.Range("N12").Value = 4160
$tasas.Range("O12").Value = 106.333
